# Fruta / hortaliza, semanal
# Inserts the latest weekly price-report row for "Caqui" (Vega Monumental
# Concepción) above the existing row 12, pushing the previous rows 12-17
# down to 13-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 12:17 down to 13:18 by inserting a new row at 12.
$ws.Rows.Item(12).Insert()

# New row 12 data (week of 2023-06-20).
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Vega Monumental Concepción"
$ws.Range("C12").Value = "Bíobío"
$ws.Range("D12").Value = 45097
$ws.Range("D12").NumberFormat = $ws.Range("D11").NumberFormat
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = "Otros"
$ws.Range("I12").Value = 100107001
$ws.Range("J12").Value = "Caqui"
$ws.Range("K12").Value = "Mankaki"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19000
$ws.Range("Q12").Value = "$/caja 18 kilos granel"
$ws.Range("R12").Value = "Región del Maule"
$ws.Range("S12").Value = 1056
$ws.Range("T12").Value = 18
